$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview": the row describing 867d9b03-...md (previously last, with
# status "Ready for handoff") has been handed back and is now sorted first;
# the other two rows shift down by one position. All three now show the
# "Handed back: in sync with en-US" status.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "867d9b03-716f-4053-b0b1-333963e0e437.md"
$wsOverview.Range("B2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("D2").Value = "2016-03-22 23:15:29"

$wsOverview.Range("A3").Value = "ffffdf981585-01e1-44ae-becf-5b182e7b97a1.md"
$wsOverview.Range("B3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("D3").Value = "2016-03-22 23:13:42"

$wsOverview.Range("A4").Value = "ffffff885f51c1-2a4e-4071-9188-dcc1200a6f7c.md"
$wsOverview.Range("B4").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C4").Value = "Handed back: in sync with en-US"
$wsOverview.Range("D4").Value = "2016-03-22 23:13:42"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/6e0cd2a2f019169a66c55b4f09b25e0a86a0ce9f/e2e/867d9b03-716f-4053-b0b1-333963e0e437.md", "", "", "867d9b03-716f-4053-b0b1-333963e0e437.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/6e0cd2a2f019169a66c55b4f09b25e0a86a0ce9f/e2e/ffffdf981585-01e1-44ae-becf-5b182e7b97a1.md", "", "", "ffffdf981585-01e1-44ae-becf-5b182e7b97a1.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/6e0cd2a2f019169a66c55b4f09b25e0a86a0ce9f/e2e/ffffff885f51c1-2a4e-4071-9188-dcc1200a6f7c.md", "", "", "ffffff885f51c1-2a4e-4071-9188-dcc1200a6f7c.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn": same row reshuffle as above, and the 867d9b03 row (now row 2)
# gets its Latest Target File / Latest Handback File / Latest Handback
# DateTime columns filled in since it has now been handed back.
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = "867d9b03-716f-4053-b0b1-333963e0e437.md"
$wsZh.Range("B2").Value = ".md"
$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Range("D2").Value = "867d9b03-716f-4053-b0b1-333963e0e437.e46e4cb8e022126f5dd0e6fea102368fb0575fba.zh-cn.xlf"
$wsZh.Range("E2").Value = "2016-03-22 23:15:22"
$wsZh.Range("F2").Value = "867d9b03-716f-4053-b0b1-333963e0e437.md"
$wsZh.Range("G2").Value = "867d9b03-716f-4053-b0b1-333963e0e437.e46e4cb8e022126f5dd0e6fea102368fb0575fba.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-03-22 23:15:55"
$wsZh.Range("J2").Value = "Include"

$wsZh.Range("A3").Value = "ffffdf981585-01e1-44ae-becf-5b182e7b97a1.md"
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"
$wsZh.Range("D3").Value = "4deb35a5-630e-4299-b40e-8f9a90586cb4.e1a71853c02e1bd31a2bcad5a76af3a7ebbee551.zh-cn.xlf"
$wsZh.Range("E3").Value = "2016-03-22 23:13:38"
$wsZh.Range("F3").Value = "4deb35a5-630e-4299-b40e-8f9a90586cb4.md"
$wsZh.Range("G3").Value = "4deb35a5-630e-4299-b40e-8f9a90586cb4.e1a71853c02e1bd31a2bcad5a76af3a7ebbee551.zh-cn.xlf"
$wsZh.Range("H3").Value = "2016-03-22 23:14:01"
$wsZh.Range("J3").Value = "Include"

$wsZh.Range("A4").Value = "ffffff885f51c1-2a4e-4071-9188-dcc1200a6f7c.md"
$wsZh.Range("B4").Value = ".md"
$wsZh.Range("C4").Value = "Handed back: in sync with en-US"
$wsZh.Range("D4").Value = "4deb35a5-630e-4299-b40e-8f9a90586cb4.e1a71853c02e1bd31a2bcad5a76af3a7ebbee551.zh-cn.xlf"
$wsZh.Range("E4").Value = "2016-03-22 23:13:38"
$wsZh.Range("F4").Value = "4deb35a5-630e-4299-b40e-8f9a90586cb4.md"
$wsZh.Range("G4").Value = "4deb35a5-630e-4299-b40e-8f9a90586cb4.e1a71853c02e1bd31a2bcad5a76af3a7ebbee551.zh-cn.xlf"
$wsZh.Range("H4").Value = "2016-03-22 23:14:01"
$wsZh.Range("J4").Value = "Include"

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/6e0cd2a2f019169a66c55b4f09b25e0a86a0ce9f/e2e/867d9b03-716f-4053-b0b1-333963e0e437.md", "", "", "867d9b03-716f-4053-b0b1-333963e0e437.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/43a403bb88fea519b8bad8c2d31e73650860fcd9/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/867d9b03-716f-4053-b0b1-333963e0e437.e46e4cb8e022126f5dd0e6fea102368fb0575fba.zh-cn.xlf", "", "", "867d9b03-716f-4053-b0b1-333963e0e437.e46e4cb8e022126f5dd0e6fea102368fb0575fba.zh-cn.xlf") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/43a403bb88fea519b8bad8c2d31e73650860fcd9/e2e/867d9b03-716f-4053-b0b1-333963e0e437.md", "", "", "867d9b03-716f-4053-b0b1-333963e0e437.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/43a403bb88fea519b8bad8c2d31e73650860fcd9/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/867d9b03-716f-4053-b0b1-333963e0e437.e46e4cb8e022126f5dd0e6fea102368fb0575fba.zh-cn.xlf", "", "", "867d9b03-716f-4053-b0b1-333963e0e437.e46e4cb8e022126f5dd0e6fea102368fb0575fba.zh-cn.xlf") | Out-Null

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/6e0cd2a2f019169a66c55b4f09b25e0a86a0ce9f/e2e/ffffdf981585-01e1-44ae-becf-5b182e7b97a1.md", "", "", "ffffdf981585-01e1-44ae-becf-5b182e7b97a1.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/02b07f23aaa38d419af93ba9023d040aa5c23598/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/4deb35a5-630e-4299-b40e-8f9a90586cb4.e1a71853c02e1bd31a2bcad5a76af3a7ebbee551.zh-cn.xlf", "", "", "4deb35a5-630e-4299-b40e-8f9a90586cb4.e1a71853c02e1bd31a2bcad5a76af3a7ebbee551.zh-cn.xlf") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/79f5ba44c16af142ad518a6816be9d29c1699b72/e2e/4deb35a5-630e-4299-b40e-8f9a90586cb4.md", "", "", "4deb35a5-630e-4299-b40e-8f9a90586cb4.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/27e98b045886f20a2eba458f53156c219ab0f2a0/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/4deb35a5-630e-4299-b40e-8f9a90586cb4.e1a71853c02e1bd31a2bcad5a76af3a7ebbee551.zh-cn.xlf", "", "", "4deb35a5-630e-4299-b40e-8f9a90586cb4.e1a71853c02e1bd31a2bcad5a76af3a7ebbee551.zh-cn.xlf") | Out-Null

$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/6e0cd2a2f019169a66c55b4f09b25e0a86a0ce9f/e2e/ffffff885f51c1-2a4e-4071-9188-dcc1200a6f7c.md", "", "", "ffffff885f51c1-2a4e-4071-9188-dcc1200a6f7c.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/02b07f23aaa38d419af93ba9023d040aa5c23598/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/4deb35a5-630e-4299-b40e-8f9a90586cb4.e1a71853c02e1bd31a2bcad5a76af3a7ebbee551.zh-cn.xlf", "", "", "4deb35a5-630e-4299-b40e-8f9a90586cb4.e1a71853c02e1bd31a2bcad5a76af3a7ebbee551.zh-cn.xlf") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("F4"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/79f5ba44c16af142ad518a6816be9d29c1699b72/e2e/4deb35a5-630e-4299-b40e-8f9a90586cb4.md", "", "", "4deb35a5-630e-4299-b40e-8f9a90586cb4.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("G4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/27e98b045886f20a2eba458f53156c219ab0f2a0/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/4deb35a5-630e-4299-b40e-8f9a90586cb4.e1a71853c02e1bd31a2bcad5a76af3a7ebbee551.zh-cn.xlf", "", "", "4deb35a5-630e-4299-b40e-8f9a90586cb4.e1a71853c02e1bd31a2bcad5a76af3a7ebbee551.zh-cn.xlf") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de": identical reshuffle/enrichment as zh-cn but with the de-de
# dated artifacts.
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = "867d9b03-716f-4053-b0b1-333963e0e437.md"
$wsDe.Range("B2").Value = ".md"
$wsDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDe.Range("D2").Value = "867d9b03-716f-4053-b0b1-333963e0e437.e46e4cb8e022126f5dd0e6fea102368fb0575fba.de-de.xlf"
$wsDe.Range("E2").Value = "2016-03-22 23:15:29"
$wsDe.Range("F2").Value = "867d9b03-716f-4053-b0b1-333963e0e437.md"
$wsDe.Range("G2").Value = "867d9b03-716f-4053-b0b1-333963e0e437.e46e4cb8e022126f5dd0e6fea102368fb0575fba.de-de.xlf"
$wsDe.Range("H2").Value = "2016-03-22 23:16:01"
$wsDe.Range("J2").Value = "Include"

$wsDe.Range("A3").Value = "ffffdf981585-01e1-44ae-becf-5b182e7b97a1.md"
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDe.Range("D3").Value = "4deb35a5-630e-4299-b40e-8f9a90586cb4.e1a71853c02e1bd31a2bcad5a76af3a7ebbee551.de-de.xlf"
$wsDe.Range("E3").Value = "2016-03-22 23:13:42"
$wsDe.Range("F3").Value = "4deb35a5-630e-4299-b40e-8f9a90586cb4.md"
$wsDe.Range("G3").Value = "4deb35a5-630e-4299-b40e-8f9a90586cb4.e1a71853c02e1bd31a2bcad5a76af3a7ebbee551.de-de.xlf"
$wsDe.Range("H3").Value = "2016-03-22 23:14:09"
$wsDe.Range("J3").Value = "Include"

$wsDe.Range("A4").Value = "ffffff885f51c1-2a4e-4071-9188-dcc1200a6f7c.md"
$wsDe.Range("B4").Value = ".md"
$wsDe.Range("C4").Value = "Handed back: in sync with en-US"
$wsDe.Range("D4").Value = "4deb35a5-630e-4299-b40e-8f9a90586cb4.e1a71853c02e1bd31a2bcad5a76af3a7ebbee551.de-de.xlf"
$wsDe.Range("E4").Value = "2016-03-22 23:13:42"
$wsDe.Range("F4").Value = "4deb35a5-630e-4299-b40e-8f9a90586cb4.md"
$wsDe.Range("G4").Value = "4deb35a5-630e-4299-b40e-8f9a90586cb4.e1a71853c02e1bd31a2bcad5a76af3a7ebbee551.de-de.xlf"
$wsDe.Range("H4").Value = "2016-03-22 23:14:09"
$wsDe.Range("J4").Value = "Include"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/6e0cd2a2f019169a66c55b4f09b25e0a86a0ce9f/e2e/867d9b03-716f-4053-b0b1-333963e0e437.md", "", "", "867d9b03-716f-4053-b0b1-333963e0e437.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/57473691255ff0419771e29064a2cd1f448ca0dd/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/867d9b03-716f-4053-b0b1-333963e0e437.e46e4cb8e022126f5dd0e6fea102368fb0575fba.de-de.xlf", "", "", "867d9b03-716f-4053-b0b1-333963e0e437.e46e4cb8e022126f5dd0e6fea102368fb0575fba.de-de.xlf") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/57473691255ff0419771e29064a2cd1f448ca0dd/e2e/867d9b03-716f-4053-b0b1-333963e0e437.md", "", "", "867d9b03-716f-4053-b0b1-333963e0e437.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/57473691255ff0419771e29064a2cd1f448ca0dd/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/867d9b03-716f-4053-b0b1-333963e0e437.e46e4cb8e022126f5dd0e6fea102368fb0575fba.de-de.xlf", "", "", "867d9b03-716f-4053-b0b1-333963e0e437.e46e4cb8e022126f5dd0e6fea102368fb0575fba.de-de.xlf") | Out-Null

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/6e0cd2a2f019169a66c55b4f09b25e0a86a0ce9f/e2e/ffffdf981585-01e1-44ae-becf-5b182e7b97a1.md", "", "", "ffffdf981585-01e1-44ae-becf-5b182e7b97a1.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/14990f17cbe64192b2773fd6fc146bbe54ec3ff2/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/4deb35a5-630e-4299-b40e-8f9a90586cb4.e1a71853c02e1bd31a2bcad5a76af3a7ebbee551.de-de.xlf", "", "", "4deb35a5-630e-4299-b40e-8f9a90586cb4.e1a71853c02e1bd31a2bcad5a76af3a7ebbee551.de-de.xlf") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/1282a3a9ff335bbd6c1b97ad549b6c344997b34a/e2e/4deb35a5-630e-4299-b40e-8f9a90586cb4.md", "", "", "4deb35a5-630e-4299-b40e-8f9a90586cb4.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/fe3ce018940fd6c70cb5600b2c8f3099a3dac2be/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/4deb35a5-630e-4299-b40e-8f9a90586cb4.e1a71853c02e1bd31a2bcad5a76af3a7ebbee551.de-de.xlf", "", "", "4deb35a5-630e-4299-b40e-8f9a90586cb4.e1a71853c02e1bd31a2bcad5a76af3a7ebbee551.de-de.xlf") | Out-Null

$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/6e0cd2a2f019169a66c55b4f09b25e0a86a0ce9f/e2e/ffffff885f51c1-2a4e-4071-9188-dcc1200a6f7c.md", "", "", "ffffff885f51c1-2a4e-4071-9188-dcc1200a6f7c.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/14990f17cbe64192b2773fd6fc146bbe54ec3ff2/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/4deb35a5-630e-4299-b40e-8f9a90586cb4.e1a71853c02e1bd31a2bcad5a76af3a7ebbee551.de-de.xlf", "", "", "4deb35a5-630e-4299-b40e-8f9a90586cb4.e1a71853c02e1bd31a2bcad5a76af3a7ebbee551.de-de.xlf") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("F4"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/1282a3a9ff335bbd6c1b97ad549b6c344997b34a/e2e/4deb35a5-630e-4299-b40e-8f9a90586cb4.md", "", "", "4deb35a5-630e-4299-b40e-8f9a90586cb4.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("G4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/fe3ce018940fd6c70cb5600b2c8f3099a3dac2be/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/4deb35a5-630e-4299-b40e-8f9a90586cb4.e1a71853c02e1bd31a2bcad5a76af3a7ebbee551.de-de.xlf", "", "", "4deb35a5-630e-4299-b40e-8f9a90586cb4.e1a71853c02e1bd31a2bcad5a76af3a7ebbee551.de-de.xlf") | Out-Null
